# fix(publipostage): Correct status name
#
# - statut_label "bleu" -> "noir"
# - statut_name "résultat et / ou publication posté" -> "résultat postés ou publiés"
# - statut_name "pas de résultat ni de publication" -> "pas de résultat postés ni publiés"
# - statut_name "résultat et / ou publication posté dans les 36 mois" -> "résultat postés ou publiés dans les 36 mois"
# - statut_name "résultat et / ou publication posté dans les 12 mois" -> "résultat postés ou publiés dans les 12 mois"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$whole = [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole

# Replace the longer/more specific phrases before the shorter one they
# contain as a prefix; using xlWhole (whole-cell match) makes this safe
# regardless of order, but it's kept this way for clarity.
$ws.Cells.Replace("résultat et / ou publication posté dans les 36 mois", "résultat postés ou publiés dans les 36 mois", $whole)
$ws.Cells.Replace("résultat et / ou publication posté dans les 12 mois", "résultat postés ou publiés dans les 12 mois", $whole)
$ws.Cells.Replace("résultat et / ou publication posté", "résultat postés ou publiés", $whole)
$ws.Cells.Replace("pas de résultat ni de publication", "pas de résultat postés ni publiés", $whole)
$ws.Cells.Replace("bleu", "noir", $whole)
